# Insert a new column at A, shifting the existing "Allowance Range"/Test 1-5
# table one column to the right, then populate the new column A with the
# "Applied Torque" values and turn the existing numeric Test columns into
# text values (matching the new "Min-Max Allowance" / Test N layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift B:F -> C:G (and old A -> B) by inserting a blank column at A.
$ws.Columns("A:A").Insert()

# Copy the header formatting (bold, border, centered) from the shifted
# header cell into the new header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("A1").Value = "Applied Torque"
$ws.Range("B1").Value = "Min-Max Allowance"
$ws.Range("C1").Value = "Test 1"
$ws.Range("D1").Value = "Test 2"
$ws.Range("E1").Value = "Test 3"
$ws.Range("F1").Value = "Test 4"
$ws.Range("G1").Value = "Test 5"

# --- Data rows: applied torque values go into the new column A as text,
# and every Test column becomes a text value (was numeric) ---
$rows = @(
    @{ Row = 2; Torque = "550"; C = "544.0"; D = "554.3"; E = "563.2"; F = "550.7"; G = "567.1" },
    @{ Row = 3; Torque = "350"; C = "355.1"; D = "363.0"; E = "355.6"; F = "359.2"; G = "349.4" },
    @{ Row = 4; Torque = "200"; C = "203.3"; D = "207.9"; E = "207.4"; F = "207.0"; G = "201.0" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Torque
    $cellA.Style = "Normal"

    foreach ($col in @(3, 4, 5, 6, 7)) {
        $letter = @{3="C"; 4="D"; 5="E"; 6="F"; 7="G"}[$col]
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $r[$letter]
        $cell.Style = "Normal"
    }
}
